$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.756.19'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '2.076.60'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.16'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.623'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('B7').Value = 'Solana'
$ws.Range('C7').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.50'
$ws.Range('E7').Value = '  +1.37%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.394'
$ws.Range('E9').Value = '  +1.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0784'
$ws.Range('E10').Value = '  +1.12%  '
$ws.Range('E11').Value = '  +3.34%  '
$ws.Range('D12').Value = '2.380.37'
$ws.Range('E12').Value = '  -1.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.77'
$ws.Range('E13').Value = '  +2.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.91'
$ws.Range('E14').Value = '  -1.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.774'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.34'
$ws.Range('E16').Value = '  +2.85%  '
$ws.Range('D17').Value = '2.109.00'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').Value = '37.682.45'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.19'
$ws.Range('E19').Value = '  -0.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.07'
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('D21').Value = '0.0₃0834'
$ws.Range('E21').Value = '  +1.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.28'
$ws.Range('E22').Value = '  +0.62%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  -2.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.40'
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '170.96'
$ws.Range('E26').Value = '  +1.34%  '
$ws.Range('E27').Value = '  +4.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.02'
$ws.Range('E28').Value = '  +0.98%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.47'
$ws.Range('E29').Value = '  +0.48%  '
$ws.Range('E30').Value = '  -2.06%  '
$ws.Range('E31').Value = '  +2.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.68'
$ws.Range('E32').Value = '  +1.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0632'
$ws.Range('E33').Value = '  +1.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.65'
$ws.Range('E34').Value = '  +2.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.49'
$ws.Range('E35').Value = '  -3.21%  '
$ws.Range('E36').Value = '  +0.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.41'
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.31'
$ws.Range('E39').Value = '  -2.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '101.51'
$ws.Range('E40').Value = '  +5.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0974'
$ws.Range('E41').Value = '  -3.00%  '
$ws.Range('E42').Value = '  -1.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0215'
$ws.Range('E43').Value = '  +1.07%  '
$ws.Range('D44').Value = '1.445.13'
$ws.Range('E44').Value = '  -0.94%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.16'
$ws.Range('E45').Value = '  -0.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.38'
$ws.Range('E46').Value = '  +7.04%  '
$ws.Range('E47').Value = '  +2.69%  '
$ws.Range('E48').Value = '  +0.87%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.43'
$ws.Range('E49').Value = '  +2.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.01'
$ws.Range('E50').Value = '  -0.78%  '
$ws.Range('D51').Value = '2.265.18'
$ws.Range('E51').Value = '  -1.52%  '
